$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, pushing existing rows 132-148 down to 133-149
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new weekly price record
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(132, 3).Value = "Ñuble"
$ws.Cells.Item(132, 4).Value = 45106
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112031
$ws.Cells.Item(132, 7).Value = "Poroto verde"
$ws.Cells.Item(132, 8).Value = "Magnum"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 40
$ws.Cells.Item(132, 11).Value = 25000
$ws.Cells.Item(132, 12).Value = 25000
$ws.Cells.Item(132, 13).Value = 25000
$ws.Cells.Item(132, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(132, 15).Value = "Perú"
$ws.Cells.Item(132, 16).Value = 1000
$ws.Cells.Item(132, 17).Value = 25
$ws.Cells.Item(132, 18).Value = "Hortaliza"
